$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared description/code for rows 2-5 (BAV24G0I1C ball valve)
$bav24Code = "BAV24G0I1C"
$bav24Desc = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"

$ws.Range("B2:B5").Value = $bav24Code
$ws.Range("C2:C5").Value = $bav24Desc

# Row 2
$ws.Range("D2").Value = "0,5"
$ws.Range("E2").Value = "1,00"
$ws.Range("F2").Value = "e.a"
$ws.Range("G2").Value = "-"

# Row 3
$ws.Range("D3").Value = "0,75"
$ws.Range("E3").Value = "1,00"
$ws.Range("F3").Value = "e.a"
$ws.Range("G3").Value = "CSO"

# Row 4
$ws.Range("D4").Value = "1,00"
$ws.Range("E4").Value = "3,00"
$ws.Range("F4").Value = "e.a"
$ws.Range("G4").Value = "-"

# Row 5
$ws.Range("D5").Value = "1,00"
$ws.Range("E5").Value = "1,00"
$ws.Range("F5").Value = "e.a"
$ws.Range("G5").Value = "CSO"

# Row 6 - swing check valve
$ws.Range("B6").Value = "CKV21A0B2B"
$ws.Range("C6").Value = "SWING CHECK VALVE FL, API 594, API 598, A216 GR.WCB, CL 150, INST HORIZ/VERT, RF, B16.5, BOLTED COVER, SPW SS304/GRAPH, RENEWABLE SEATS, TRIM #8"
$ws.Range("D6").Value = "3,00"
$ws.Range("E6").Value = "1,00"
$ws.Range("F6").Value = "e.a"
$ws.Range("G6").Value = "-"

# Row 7 - lift check valve
$ws.Range("B7").Value = "CLV24F0B2B"
$ws.Range("C7").Value = "LIFT CHECK VALVE SW, API 602, A105, CL 800, INST HORIZ/VERT, SW, B16.11, BOLTED COVER, SPW SS304/GRAPH, PISTON TYPE OBTURATOR, TRIM #8"
$ws.Range("D7").Value = "0,5"
$ws.Range("E7").Value = "1,00"
$ws.Range("F7").Value = "e.a"
$ws.Range("G7").Value = "-"

# Row 8 - gate valve
$ws.Range("B8").Value = "GAV24F0B2B"
$ws.Range("C8").Value = "GATE VALVE SW, API 602, API 598, A105, CL 800, SW, B16.11, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, SOLID WEDGE, STEM OS&Y/RSNRO, HO"
$ws.Range("D8").Value = "0,75"
$ws.Range("E8").Value = "2,00"
$ws.Range("F8").Value = "e.a"
$ws.Range("G8").Value = "-"

# Row 9 - integral mono flange DBB needle multi-valve
$ws.Range("B9").Value = "MFV21A0I2I"
$ws.Range("C9").Value = "INTEGRAL MONO FLANGE DBB NEEDLE MULTI-VALVE, EEMUA 182, A105, CL 150, RF/NPTF, B16.5 AND B1.20.1, BB, SPW SS304/GRAPH, PKG GRAPH; SS316 STEM, SEATS&STEM TIP, S, SWIVEL NEEDLE, STEM OS&Y/RSRO, T-HANDLE"
$ws.Range("D9").Value = "0,75"
$ws.Range("E9").Value = "6,00"
$ws.Range("F9").Value = "e.a"
$ws.Range("G9").Value = "-"

# Remove the now-obsolete rows 10-17 (CKV/GAV/GLV variants that were folded
# into the rows above or dropped entirely)
$ws.Rows("10:17").Delete()
